{"js": "const doc = context.document;\nconst body = doc.body;\n\n// 1) \"... November 2013 \u2013 October 2015\" -> \"... November 2013 \u2013 March 2015\"\nlet results = body.search(\"October 2015\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\nresults.items[0].insertText(\"March 2015\", \"Replace\");\nawait context.sync();\n\n// 2) \"...website and iOS app.\" - no visible text changes here, but the run\n//    that used to be split around the stray \"_GoBack\" bookmark (\" iOS \" /\n//    \"app\") is rejoined into a single run now that the bookmark is gone.\n//    The bookmark itself is relocated (removed here, re-inserted near the\n//    date fix below) just like Word does when the last edit position moves\n//    elsewhere in the doc.\nresults = body.search(\"er cards integration on the website and iOS app\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\nresults.items[0].insertText(\"er cards integration on the website and iOS app\", \"Replace\");\nawait context.sync();\n\n// 3) \"October 2103\" -> \"October 2013\" (typo fix)\nresults = body.search(\"October 2103\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\nresults.items[0].insertText(\"October 2013\", \"Replace\");\nawait context.sync();\n\n// 4) \"Feb 1994 \u2013 Jun 2007\" -> \"February 1994 \u2013 June 2007\"\nresults = body.search(\"Feb 1994\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\nresults.items[0].insertText(\"February 1994\", \"Replace\");\nawait context.sync();\n\nresults = body.search(\"Jun 2007\", { matchCase: true, matchWholeWord: true });\nresults.load(\"items\");\nawait context.sync();\nresults.items[0].insertText(\"June 2007\", \"Replace\");\nawait context.sync();\n\n// Move the \"_GoBack\" bookmark from the \"iOS app\" paragraph to sit right\n// after the newly expanded \"February\", mirroring Word's \"last edit\"\n// bookmark behavior.\ndoc.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nresults = body.search(\"February\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\nresults.items[0].insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) \"... November 2013 - October 2015\" -> \"... November 2013 - March 2015\"\n$r1 = $d.Range(0, 0)\n$r1.Find.Execute(\"October 2015\", $false, $true, $false, $false, $false, $true, 1, $false, \"March 2015\", 2) | Out-Null\n\n# 2) \"...website and iOS app.\" - no visible text changes here; the stray\n#    \"_GoBack\" bookmark that used to sit between \" iOS \" and \"app\" is\n#    relocated (removed here, re-inserted near the date fix below), the\n#    same way Word moves it to track the last edited spot.\n\n# 3) \"October 2103\" -> \"October 2013\" (typo fix)\n$r3 = $d.Range(0, 0)\n$r3.Find.Execute(\"October 2103\", $false, $true, $false, $false, $false, $true, 1, $false, \"October 2013\", 2) | Out-Null\n\n# 4) \"Feb 1994 - Jun 2007\" -> \"February 1994 - June 2007\"\n$r4a = $d.Range(0, 0)\n$r4a.Find.Execute(\"Feb 1994\", $false, $true, $false, $false, $false, $true, 1, $false, \"February 1994\", 2) | Out-Null\n\n$r4b = $d.Range(0, 0)\n$r4b.Find.Execute(\"Jun 2007\", $false, $true, $false, $false, $false, $true, 1, $false, \"June 2007\", 2) | Out-Null\n\n# Move the \"_GoBack\" bookmark from the \"iOS app\" paragraph to sit right\n# after the newly expanded \"February\", mirroring Word's \"last edit\"\n# bookmark behavior.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$rFeb = $d.Range(0, 0)\n$rFeb.Find.Execute(\"February\") | Out-Null\n$d.Bookmarks.Add(\"_GoBack\", $rFeb) | Out-Null\n"}
